# Update to version 6.0.0 of cht-conf
# - settings sheet: drop the "form_id" column (B) entirely - forms are no
#   longer required to declare a form_id; version/style/namespaces shift
#   left by one column and their header comments shift accordingly.
# - survey sheet: conditional-formatting ranges get re-normalised.
# - the "settings" tab becomes the active tab/sheet instead of "survey".

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------
# settings sheet: remove the form_id column (old column B).
# ---------------------------------------------------------------------

# Grab the existing header comments (for B1/C1/D1/E1) before the column
# shift so we can re-attach the right text to the right header afterwards.
$oldVersionComment     = $settings.Range("C1").Comment.Text()
$oldStyleComment       = $settings.Range("D1").Comment.Text()
$oldNamespacesComment  = $settings.Range("E1").Comment.Text()

# Delete the form_id column - this shifts version/style/namespaces left,
# fixes up the row data, dimension and shared strings automatically.
$settings.Columns.Item(2).Delete()

# Re-point each remaining header's comment at the text that used to
# belong one column to the right of it.
$settings.Range("B1").Comment.Text($oldVersionComment)
$settings.Range("C1").Comment.Text($oldStyleComment)
$settings.Range("D1").Comment.Text($oldNamespacesComment)

# The old E1 comment (namespaces) is now a duplicate of D1 - drop it.
$leftoverComment = $settings.Range("E1").Comment
if ($leftoverComment) {
    $leftoverComment.Delete()
}

# ---------------------------------------------------------------------
# survey sheet: re-normalise the conditional-formatting ranges.
# ---------------------------------------------------------------------

$mainRuleGroup = $survey.Range("A2").FormatConditions.Item(1)
$mainRuleGroup.ModifyAppliesToRange($survey.Range("A2:F10000"))

$colCRules = $survey.Range("C2").FormatConditions
for ($i = 1; $i -le $colCRules.Count(); $i++) {
    $rule = $colCRules.Item($i)
    if ($rule.Formula1() -like "*ISBLANK(C2)*") {
        $rule.ModifyAppliesToRange($survey.Range("C2:C10000"))
    }
}

# ---------------------------------------------------------------------
# make "settings" the active sheet/tab, with B10 selected on it.
# ---------------------------------------------------------------------

$settings.Activate()
$settings.Range("B10").Select()
